# Montecreto.xlsx update ("aggiornato a 2/3, aggiornati i report")
#
# Effect (per the OOXML diff):
#  - A new daily data point for 2021-02-08 (serial 44235) is inserted
#    between the existing 2021-02-07 (44234) and 2021-02-09 (44236) rows,
#    shifting every following row down by one.
#  - The rolling "somma mobile 7gg." (col C) and the per-100k-inhabitants
#    figure (col D) are refreshed for the affected window.
#  - Two more days of (still-empty-C/D) data are appended at the bottom
#    (2021-03-01 / 44256 and 2021-03-02 / 44257), extending the sheet
#    from A1:D113 to A1:D115.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new row for 2021-02-08 (44235) right after row 92 (44234) ---
# This pushes the old rows 93..113 down to 94..114, carrying their existing
# values/formatting with them.
$ws.Rows("93:93").Insert()

# The freshly-inserted row 93 has no content yet, and the brand-new row 115
# (appended past the old end of the table) doesn't exist yet either. Give
# both of their date cells (column A) the same look as the rest of column A
# by copying the format from an existing, correctly-styled date cell.
$ws.Range("A92").Copy()
$ws.Range("A93").PasteSpecial(-4122)   # -4122 = xlPasteFormats
$ws.Range("A115").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rewrite A/B/C/D for every row from the insertion point through the
#     (now extended) end of the table, using the post-edit values. ---

$rows = @(
    @{ Row = 92;  Date = 44234; B = 3; C = 4;    D = 436.6812227074236 },
    @{ Row = 93;  Date = 44235; B = 0; C = 5;    D = 545.8515283842795 },
    @{ Row = 94;  Date = 44236; B = 0; C = 6;    D = 655.0218340611353 },
    @{ Row = 95;  Date = 44237; B = 0; C = 6;    D = 655.0218340611353 },
    @{ Row = 96;  Date = 44238; B = 2; C = 5;    D = 545.8515283842795 },
    @{ Row = 97;  Date = 44239; B = 1; C = 5;    D = 545.8515283842795 },
    @{ Row = 98;  Date = 44240; B = 0; C = 5;    D = 545.8515283842795 },
    @{ Row = 99;  Date = 44241; B = 2; C = 5;    D = 545.8515283842795 },
    @{ Row = 100; Date = 44242; B = 0; C = 4;    D = 436.6812227074236 },
    @{ Row = 101; Date = 44243; B = 0; C = 3;    D = 327.5109170305677 },
    @{ Row = 102; Date = 44244; B = 0; C = 3;    D = 327.5109170305677 },
    @{ Row = 103; Date = 44245; B = 1; C = 2;    D = 218.3406113537118 },
    @{ Row = 104; Date = 44246; B = 0; C = 2;    D = 218.3406113537118 },
    @{ Row = 105; Date = 44247; B = 0; C = 4;    D = 436.6812227074236 },
    @{ Row = 106; Date = 44248; B = 1; C = 4;    D = 436.6812227074236 },
    @{ Row = 107; Date = 44249; B = 0; C = 4;    D = 436.6812227074236 },
    @{ Row = 108; Date = 44250; B = 2; C = 4;    D = 436.6812227074236 },
    @{ Row = 109; Date = 44251; B = 0; C = 4;    D = 436.6812227074236 },
    @{ Row = 110; Date = 44252; B = 1; C = 3;    D = 327.5109170305677 },
    @{ Row = 111; Date = 44253; B = 0; C = 3;    D = 327.5109170305677 },
    @{ Row = 112; Date = 44254; B = 0; C = 1;    D = 109.1703056768559 },
    @{ Row = 113; Date = 44255; B = 0; C = $null; D = $null },
    @{ Row = 114; Date = 44256; B = 0; C = $null; D = $null },
    @{ Row = 115; Date = 44257; B = 0; C = $null; D = $null }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.B

    if ($null -ne $r.C) {
        $ws.Cells.Item($r.Row, 3).Value = $r.C
        $ws.Cells.Item($r.Row, 4).Value = $r.D
    }
}

# The trailing rows (113-115) keep columns C/D blank, exactly like the other
# not-yet-reported days elsewhere in the sheet (e.g. C2/D2): present as
# empty-string cells rather than fully absent ones, matching the workbook's
# existing convention for "no data yet".
$ws.Range("C2:D2").Copy()
$ws.Range("C113:D115").PasteSpecial(-4163)   # -4163 = xlPasteValues
$excel.CutCopyMode = $false
